{"js": "// Remove the \"Instructions\" bullet list (dataset/codes/turn-in/submit/\n// deadline items), the trailing blank list item, the \"Questions\" heading,\n// and the blank paragraph that followed it \u2014 i.e. everything between the\n// \"Instructions\" heading paragraph and the \"The Boston Housing data ...\"\n// paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst startMarker = \"The dataset\";\nconst endMarker = \"The Boston Housing data\";\n\nlet deleting = false;\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text || \"\";\n  if (!deleting && text.indexOf(startMarker) !== -1) {\n    deleting = true;\n  }\n  if (deleting && text.indexOf(endMarker) !== -1) {\n    break;\n  }\n  if (deleting) {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Instructions\" bullet list (dataset/codes/turn-in/submit/\n# deadline items), the trailing blank list item, the \"Questions\" heading,\n# and the blank paragraph that followed it -- i.e. everything between the\n# \"Instructions\" heading paragraph and the \"The Boston Housing data ...\"\n# paragraph.\n$d = $word.ActiveDocument\n\n$startRange = $d.Content\n$startRange.Find.Execute(\"The dataset\")\n$startPara = $startRange.Paragraphs.Item(1)\n$start = $startPara.Range.Start\n\n$endRange = $d.Content\n$endRange.Find.Execute(\"The Boston Housing data\")\n$endPara = $endRange.Paragraphs.Item(1)\n$end = $endPara.Range.Start\n\n$deleteRange = $d.Range($start, $end)\n$deleteRange.Delete()\n"}
